$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 65, shifting existing rows 65-98 down to 66-99.
$ws.Rows.Item(65).Insert()

# Populate the new row 65 with the new record (Florida King, Región de Coquimbo).
$ws.Cells.Item(65, 1).Value = 1
$ws.Cells.Item(65, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(65, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(65, 4).Value = 45240
$ws.Cells.Item(65, 5).Value = 15
$ws.Cells.Item(65, 6).Value = "Fruta"
$ws.Cells.Item(65, 7).Value = 100103
$ws.Cells.Item(65, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(65, 9).Value = 100103004
$ws.Cells.Item(65, 10).Value = "Durazno"
$ws.Cells.Item(65, 11).Value = "Florida King"
$ws.Cells.Item(65, 12).Value = "Segunda"
$ws.Cells.Item(65, 13).Value = 250
$ws.Cells.Item(65, 14).Value = 31000
$ws.Cells.Item(65, 15).Value = 33000
$ws.Cells.Item(65, 16).Value = 32000
$ws.Cells.Item(65, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(65, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(65, 19).Value = 1778
$ws.Cells.Item(65, 20).Value = 18
